$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.103.13'
$ws.Range("E2").Value = '  +3.25%  '
$ws.Range("D3").Value = '1.656.51'
$ws.Range("E3").Value = '  +3.75%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.26'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.79%  '
$ws.Range("E6").Value = '  +0.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.13%  '
$ws.Range("E8").Value = '  +1.72%  '
$ws.Range("E9").Value = '  +1.54%  '
$ws.Range("E10").Value = '  +3.56%  '
$ws.Range("E11").Value = '  +1.30%  '
$ws.Range("D12").Value = '1.890.14'
$ws.Range("E12").Value = '  +3.70%  '
$ws.Range("D13").Value = '1.662.58'
$ws.Range("E13").Value = '  +4.16%  '
$ws.Range("E14").Value = '  +2.12%  '
$ws.Range("E15").Value = '  +3.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.00'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.15%  '
$ws.Range("D17").Value = '27.069.29'
$ws.Range("E17").Value = '  +3.10%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '237.98'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.97'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.06%  '
$ws.Range("D20").Value = '0.0₃0730'
$ws.Range("E20").Value = '  +1.14%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.44'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +4.90%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.24'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.61%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '9.26'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.61%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '145.92'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.16'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.05%  '
$ws.Range("E28").Value = '  +1.20%  '
$ws.Range("E29").Value = '  +3.45%  '
$ws.Range("E30").Value = '  +0.79%  '
$ws.Range("E31").Value = '  +1.65%  '
$ws.Range("D32").Value = '1.527.33'
$ws.Range("E32").Value = '  +4.28%  '
$ws.Range("E33").Value = '  +2.83%  '
$ws.Range("E34").Value = '  +3.81%  '
$ws.Range("E35").Value = '  +8.32%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.42'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.08%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.578'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.91%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.886'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.10%  '
$ws.Range("E39").Value = '  +2.75%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '5.95'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.02%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '66.43'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.76%  '
$ws.Range("E43").Value = '  +3.78%  '
$ws.Range("D44").Value = '1.796.96'
$ws.Range("E44").Value = '  +3.51%  '
$ws.Range("E45").Value = '  +2.78%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '90.20'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.83%  '
$ws.Range("D48").Value = '0.0₆0105'
$ws.Range("E48").Value = '  +0.67%  '
$ws.Range("E49").Value = '  +3.28%  '
$ws.Range("E50").Value = '  +0.84%  '
$ws.Range("E51").Value = '  +3.21%  '
